$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.304.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.500.11'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.44'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.09'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.80%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.26'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.08%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.092.56'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000182'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.498.50'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.82'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.287.75'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.89'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.75'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.62'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '393.56'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.569'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.639.39'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.36'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.60'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000115'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.41'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.50'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.37%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.24'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.99%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.27'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.519.22'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.150'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.63%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.42'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.16'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.90'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.44'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.86%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0782'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.806'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.26'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.88%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.40'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.18'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.58%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.65'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.97%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.460.31'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.77'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.896'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0261'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.34%  '
